$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at 12 (shifts existing rows 12-29 down to 13-30)
$ws.Rows("12:12").Insert()

# 2. Restore the row height for the new row 12 (Insert leaves it blank/default)
$ws.Rows("12:12").RowHeight = 25.5

# 3. Copy cell formatting from row 13 (the row right below, which holds what
#    used to be row 12's formatting) onto the new row 12
$ws.Range("A13:Q13").Copy()
$ws.Range("A12:Q12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Re-create the merged cells for row 12 to match the rest of the table
$ws.Range("A12:B12").Merge()
$ws.Range("C12:G12").Merge()
$ws.Range("H12:K12").Merge()
$ws.Range("L12:M12").Merge()
$ws.Range("N12:O12").Merge()

# 5. Fill in the values for the new item (#6 - DEPOVIT)
$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "DEPOVIT B12-1000MCG/ML 5 I.M. AMP"
$ws.Range("H12").Value = "'2:3"
$ws.Range("L12").Value = "'1"
$ws.Range("N12").Value = "'85.00"
$ws.Range("P12").Value = "'85.0000"
$ws.Range("Q12").Value = "'1:0"

# 6. Update "سرنجات 3 سم" row (now row 26) sale price / transaction count
$ws.Range("P26").Value = "'18.0000"
$ws.Range("Q26").Value = "'9:0"

# 7. Update the grand total (now row 29)
$ws.Range("P29").Value = 1198.59

# 8. Update the printed timestamp (now row 30)
$ws.Range("A30").Value = "Thursday, 29 May, 2025 1:55 PM"

Write-Output "done"
